$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timestamp value pushed into the top block (rows 2-15)
$newTop = 44262.50961724984

# Values previously held by the block above shift down into the next block
$midVal = 44262.4882553588     # new value for rows 16-29 (was 44262.46689984954)
$lowVal = 44262.46689984954    # new value for rows 30-43 (was 44262.4455540625)

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value = $newTop
}

for ($r = 16; $r -le 29; $r++) {
    $ws.Cells.Item($r, 4).Value = $midVal
}

for ($r = 30; $r -le 43; $r++) {
    $ws.Cells.Item($r, 4).Value = $lowVal
}
